# Add season record columns (Wins, Losses, Ties) to the HOU_2011 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column labels in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, bordered, centered) already used by
# the rest of row 1, by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-49): same season record repeated on every player row ---
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 56
    $ws.Cells.Item($r, 31).Value = 106
    $ws.Cells.Item($r, 32).Value = 0
}
